# EBEWE Dates.xlsx - account for new EBEWE exemption changes
#
# - "COMPLIANCE DUE DATE" column (B): the 2021/2022 compliance dates are no
#   longer marked with the "*" footnote qualifier.
# - "ENERGY CONSUMPTION COMPARATIVE PERIOD" column (D): the old fixed
#   5-year windows ("Jan 20xx - Dec 20yy") are replaced with rolling
#   "Data must not be older than 12/1/20xx" wording.
# - Refresh the sheet selection / print orientation left behind by the
#   author's last save.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- COMPLIANCE DUE DATE column: drop the trailing "*" footnote marker ---
$ws.Range("B2").Value = "Dec 1, 2021"
$ws.Range("B3").Value = "Dec 1, 2022"

# --- ENERGY CONSUMPTION COMPARATIVE PERIOD column: new wording ---
$ws.Range("D2").Value = "Data must not be older than 12/1/2016"
$ws.Range("D3").Value = "Data must not be older than 12/1/2017"
$ws.Range("D4").Value = "Data must not be older than 12/1/2018"
$ws.Range("D5").Value = "Data must not be older than 12/1/2019"
$ws.Range("D6").Value = "Data must not be older than 12/1/2020"

# --- misc view/print state left over from the author's last save ---
$ws.Range("D7").Select() | Out-Null
$ws.PageSetup.Orientation = 1
